# Automatic update of files.
# Two species-observation records (row 2 and row 3) had their
# record-specific fields swapped between the rows while the shared
# location/date/observer fields stayed on the same row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values are exchanged between row 2 and row 3:
# A  Id
# B  Taxonsorteringsordning
# E  TaxonId
# F  Artnamn
# G  Vetenskapligt namn
# H  Auktor
# I  Antal
# Q  Ost
# R  Nord
$columns = @("A", "B", "E", "F", "G", "H", "I", "Q", "R")

foreach ($col in $columns) {
    $cellRow2 = $ws.Range($col + "2")
    $cellRow3 = $ws.Range($col + "3")

    $valueRow2 = $cellRow2.Value()
    $valueRow3 = $cellRow3.Value()

    $cellRow2.Value = $valueRow3
    $cellRow3.Value = $valueRow2
}
